$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill new investor rows (5-8), column by column to match the original ---
# --- authoring order of the shared-string table.                          ---

# Column A - Investor
$ws.Range("A5").Value = "Investor 3"
$ws.Range("A6").Value = "Investor 4"
$ws.Range("A7").Value = "Investor 5"
$ws.Range("A8").Value = "Investor 6"

# Column B - First Name
$ws.Range("B5").Value = "Emp3"
$ws.Range("B6").Value = "Emp4"
$ws.Range("B7").Value = "Emp5"
$ws.Range("B8").Value = "Emp6"

# Column C - Last Name
$ws.Range("C5").Value = "L3"
$ws.Range("C6").Value = "L4"
$ws.Range("C7").Value = "L5"
$ws.Range("C8").Value = "L6"

# Column D - Email (with mailto hyperlinks)
$ws.Range("D5").Value = "emp1@investor3.com"
$ws.Range("D6").Value = "emp1@investor4.com"
$ws.Range("D7").Value = "emp1@investor5.com"
$ws.Range("D8").Value = "emp1@investor6.com"

# Column E - Email Enabled
$ws.Range("E5").Value = "Yes"
$ws.Range("E6").Value = "Yes"
$ws.Range("E7").Value = "Yes"
$ws.Range("E8").Value = "Yes"

# Column G - Country Code
$ws.Range("G5").Value = "IN(91)"
$ws.Range("G6").Value = "IN(91)"
$ws.Range("G7").Value = "IN(91)"
$ws.Range("G8").Value = "IN(91)"

# Column H - Phone (only rows 5-7 have a phone number)
$ws.Range("H5").Value = 1111111111
$ws.Range("H6").Value = 2222222222
$ws.Range("H7").Value = 3333333333

# Column I - WhatsApp Enabled
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

# Column J - Approved
$ws.Range("J5").Value = "Yes"
$ws.Range("J6").Value = "Yes"
$ws.Range("J7").Value = "Yes"

# Column K - Send Confirmation Email
$ws.Range("K5").Value = "No"
$ws.Range("K6").Value = "No"
$ws.Range("K7").Value = "No"

# --- Hyperlinks for the new emails ---
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:emp1@investor3.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:emp1@investor4.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:emp1@investor5.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:emp1@investor6.com") | Out-Null

# Restore the plain "Hyperlink" cell style (Hyperlinks.Add() above tends to
# re-apply an explicit font override); this keeps the new hyperlink cells
# consistent with the existing D2:D4 hyperlink cells.
$ws.Range("D5:D8").Style = "Hyperlink"

# --- Update the active selection to match the new extent of the sheet ---
$ws.Range("K8").Select() | Out-Null
